$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("Z2").Value = "Term of Product and Riders"
$ws.Range("AA2").Value = "Policy and Rider Term"
